$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.655.81"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.840.77"
$ws.Range("E3").Value = "  +2.88%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "359.81"
$ws.Range("E5").Value = "  +8.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "116.84"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  +3.10%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.06"
$ws.Range("E10").Value = "  +2.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0866"
$ws.Range("E11").Value = "  +4.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.10"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("E14").Value = "  +3.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.288.76"
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.828.28"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.591.11"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.35"
$ws.Range("E19").Value = "  +6.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.19"
$ws.Range("E20").Value = "  +6.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.78"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0988"
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "272.87"
$ws.Range("E23").Value = "  -1.86%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.78"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("E25").Value = "  +6.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.37"
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  +1.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.141"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.84"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.28"
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0452"
$ws.Range("E33").Value = "  +30.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.14"
$ws.Range("E36").Value = "  +1.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("E38").Value = "  +1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.94"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.98"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.64"
$ws.Range("E41").Value = "  +10.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.96"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.22"
$ws.Range("E44").Value = "  -3.22%  "
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.39"
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.073.33"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.974"
$ws.Range("E49").Value = "  +10.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.66"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.17"
$ws.Range("E51").Value = "  +1.87%  "
